# Word COM-interop script implementing the diff:
#  - adds a new changelog entry (18-9-2025) + spacer paragraph before the
#    existing "02-01-2024" entry
#  - appends " (no split)" to the "Transceivers which are checked ..." line
#  - appends " and it supports Split mode." to the "Both solutions work
#    fine ..." line
#  - cleans up a handful of spell-check-split runs (TxReq/ino/leds/
#    set_tune_pwr/...) by re-typing the surrounding text in place; this
#    merges the runs and drops the now-stale <w:proofErr/> markers without
#    altering the visible text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert the new 18-9-2025 changelog paragraph + blank spacer before
#    the existing first changelog line ("02-01-2024: ...").
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("02-01-2024", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $find.Parent
    $anchorPara = $target.Paragraphs(1)
    $insPoint = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)
    $newText = "18-9-2025: Added support for Tuning in Split Operation (VFO-B and VFO B). Originally only tuning VFO A was supported. ONLY for FM-N solution and ONLY for FTdx101D and FTdx101MP.`r`r"
    $insPoint.InsertBefore($newText)
}

# ---------------------------------------------------------------------
# 2. " (no split)" appended after the FTdx-10 / FT-991 sentence.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Transceivers which are checked by their CAT commands to be supported:   FTdx-10, FT-991", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r = $find.Parent
    $insPoint = $d.Range($r.End, $r.End)
    $insPoint.InsertAfter(" (no split)")
}

# ---------------------------------------------------------------------
# 3. Re-type a few phrases in place so Word drops the stale <w:proofErr/>
#    spell-check markers and merges the surrounding runs (no visible text
#    change).
# ---------------------------------------------------------------------
function Retype-Text($needle) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}

Retype-Text "the .ino file. And Mode FM- is available in the TXCR (unless TxREq is used)."
Retype-Text "The tune signal (TxReq) is available on pin 11"
Retype-Text "the variable; set_tune_pwr "
Retype-Text "set_tune_pwr `"PC020;`""
Retype-Text "The leds indicate if there is communication"
Retype-Text "The first uses the TxReq signal on the linear connector"
Retype-Text "Rx, Tx and Gnd) and a "
Retype-Text "2x2-wire cable (TxReq + GND and"
Retype-Text "different connector/pin for the TxReq signal."
Retype-Text "The second solution does not use the TxReq signal."
Retype-Text "(. ino file) below"
Retype-Text "navigate to the .ino file and open it."

# ---------------------------------------------------------------------
# 4. "Both solutions work fine. FM-N is recommended, it easiest." gains
#    a trailing clause.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Both solutions work fine. FM-N is recommended, it easiest.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r = $find.Parent
    # drop the trailing period, then add the new clause as its own run
    $periodRange = $d.Range($r.End - 1, $r.End)
    $periodRange.Text = ""
    $insPoint = $d.Range($r.End - 1, $r.End - 1)
    $insPoint.InsertAfter(" and it supports Split mode.")
}

Write-Output "done"
